# Clean up folder structure: remove the duplicate "Helgolander Bucht / Nordsoen /
# Sorlige Nordsjo I / Sorlige Nordsjo II" blocks (rows 623-626 and 641-644) from the
# GeneratorsOfNode sheet. These were stray duplicate entries (without the proper
# Norwegian characters) that already exist further down the list (Helgoländer Bucht /
# Nordsøen, rows 689-690 & 720-721). Removing them shifts all following rows up by 8,
# bringing the used range from A1:B741 down to A1:B733.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GeneratorsOfNode")

# Delete the second (floating) block first so the first block's row numbers
# (623-626) stay valid while we work.
$ws.Range("A641:B644").EntireRow.Delete() | Out-Null
$ws.Range("A623:B626").EntireRow.Delete() | Out-Null
